$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BOM rows: mini hooks (E-Z-Hook) ---
# Fill text columns first (order chosen to reproduce shared-string layout)
$ws.Range("A5").Value = "Mini Hook Black"
$ws.Range("A6").Value = "Mini Hook Red"
$ws.Range("B5").Value = "E-Z-Hook"
$ws.Range("B6").Value = "E-Z-Hook"
$ws.Range("C6").Value = "XR25RED"
$ws.Range("C5").Value = "XR25BLK"

# Pre-seed the Octopart link text so the shared string order matches
# (D5 = black part link, D6 = red part link)
$ws.Range("D5").Value = "https://octopart.com/xr25blk-e-z-hook-19790456"
$ws.Range("D6").Value = "https://octopart.com/xr25red-e-z-hook-19790462"

# Add the actual hyperlinks (D6 first, then D5, to match relationship order)
$ws.Hyperlinks.Add($ws.Range("D6"), "https://octopart.com/xr25red-e-z-hook-19790462")
$ws.Hyperlinks.Add($ws.Range("D5"), "https://octopart.com/xr25blk-e-z-hook-19790456")

# Match the "Hyperlink" cell style used by the other Octopart Link cells
$ws.Range("D5:D6").Style = $ws.Range("D2").Style

# Move the active selection like the author's saved state
$ws.Range("D10").Select() | Out-Null
